$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.937.83"
$ws.Range("E2").Value = "  -4.37%  "

$ws.Range("D3").Value = "2.988.66"
$ws.Range("E3").Value = "  -5.49%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.63"
$ws.Range("E5").Value = "  -3.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.62"
$ws.Range("E6").Value = "  -8.13%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.565"
$ws.Range("E8").Value = "  -4.14%  "

$ws.Range("D9").Value = "2.984.08"
$ws.Range("E9").Value = "  -5.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.111"
$ws.Range("E10").Value = "  -5.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.22"
$ws.Range("E11").Value = "  -6.99%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.365"
$ws.Range("E12").Value = "  -5.15%  "

$ws.Range("D13").Value = "3.514.58"
$ws.Range("E13").Value = "  -5.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.123"
$ws.Range("E14").Value = "  -3.73%  "

$ws.Range("D15").Value = "62.146.42"
$ws.Range("E15").Value = "  -3.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.62"
$ws.Range("E16").Value = "  -6.81%  "

$ws.Range("D17").Value = "2.996.28"
$ws.Range("E17").Value = "  -5.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000147"
$ws.Range("E18").Value = "  -5.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "390.10"
$ws.Range("E19").Value = "  -5.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.09"
$ws.Range("E20").Value = "  -3.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.91"
$ws.Range("E21").Value = "  -5.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.64"
$ws.Range("E22").Value = "  -6.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.78"
$ws.Range("E24").Value = "  -4.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.465"
$ws.Range("E25").Value = "  -4.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.184"
$ws.Range("E26").Value = "  -8.24%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.40%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0945"
$ws.Range("E28").Value = "  -8.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.43"
$ws.Range("E29").Value = "  -5.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.73"
$ws.Range("E31").Value = "  -4.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.39"
$ws.Range("E32").Value = "  -4.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.71"
$ws.Range("E33").Value = "  +1.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.64"
$ws.Range("E34").Value = "  -5.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.97"
$ws.Range("E35").Value = "  -5.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.07"
$ws.Range("E36").Value = "  -5.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("E37").Value = "  -5.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.55"
$ws.Range("E38").Value = "  -7.96%  "

$ws.Range("D39").Value = "2.440.55"
$ws.Range("E39").Value = "  -10.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.55"
$ws.Range("E40").Value = "  -3.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.89"
$ws.Range("E41").Value = "  -4.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.21"
$ws.Range("E42").Value = "  -5.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.658"
$ws.Range("E43").Value = "  -6.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0592"
$ws.Range("E44").Value = "  -6.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0245"
$ws.Range("E46").Value = "  -6.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.89"
$ws.Range("E47").Value = "  -11.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0947"
$ws.Range("E48").Value = "  -3.30%  "

$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.50"
$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.62"
$ws.Range("E50").Value = "  -7.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "261.08"
$ws.Range("E51").Value = "  -9.33%  "
